$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values from the coinranking.com refresh -----------------
# Maps cell reference -> new text. Price/volume columns (D, E) contain
# free-form text such as "27.715.36", "0.0870" or "  +4.23%  "; several
# of those look like plain numbers (e.g. "0.976") and Excel would happily
# reinterpret them as numeric values (losing the trailing zero / exact
# formatting) unless the cell is explicitly marked as text first.
$newValues = [ordered]@{
    'D2' = '28.260.85'
    'E2' = '  +4.23%  '
    'D3' = '1.576.14'
    'E3' = '  +0.39%  '
    'D4' = '0.976'
    'E4' = '  -3.62%  '
    'D5' = '210.24'
    'E5' = '  -0.28%  '
    'D6' = '0.497'
    'E6' = '  +1.06%  '
    'D7' = '0.968'
    'E7' = '  -3.85%  '
    'D8' = '23.55'
    'E8' = '  +6.17%  '
    'E9' = '  +0.66%  '
    'E10' = '  +0.06%  '
    'D11' = '0.0870'
    'E11' = '  +0.95%  '
    'D12' = '1.808.67'
    'E12' = '  +0.87%  '
    'D13' = '1.571.45'
    'E13' = '  +0.15%  '
    'E14' = '  -0.02%  '
    'D15' = '0.525'
    'E15' = '  +0.90%  '
    'D16' = '28.322.75'
    'E16' = '  +4.44%  '
    'D17' = '63.59'
    'E17' = '  +2.44%  '
    'D18' = '237.98'
    'E18' = '  +10.09%  '
    'D19' = '7.57'
    'E19' = '  +2.24%  '
    'D20' = '0.0₃0707'
    'E20' = '  +0.38%  '
    'D21' = '0.982'
    'E21' = '  -2.58%  '
    'D22' = '4.12'
    'E22' = '  -0.65%  '
    'D23' = '9.47'
    'E23' = '  +2.74%  '
    'D24' = '1.94'
    'E24' = '  -0.44%  '
    'D25' = '148.60'
    'E25' = '  -3.64%  '
    'D26' = '15.45'
    'E26' = '  +2.17%  '
    'B27' = 'Stellar'
    'C27' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D27' = '0.108'
    'E27' = '  +1.87%  '
    'B28' = 'Cosmos'
    'C28' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D28' = '6.58'
    'E28' = '  -0.56%  '
    'D29' = '0.976'
    'E29' = '  -2.99%  '
    'E30' = '  -0.43%  '
    'E31' = '  +0.18%  '
    'E32' = '  +0.26%  '
    'B33' = 'InternetComputer(DFINITY)'
    'C33' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D33' = '3.15'
    'E33' = '  -1.23%  '
    'B34' = 'Maker'
    'C34' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D34' = '1.406.94'
    'E34' = '  -1.94%  '
    'D35' = '1.61'
    'E35' = '  +0.23%  '
    'D36' = '1.06'
    'E36' = '  -4.15%  '
    'D37' = '2.30'
    'E37' = '  -2.21%  '
    'D39' = '0.550'
    'E39' = '  +3.21%  '
    'D40' = '2.44'
    'E40' = '  +2.49%  '
    'D41' = '0.818'
    'E41' = '  +0.93%  '
    'D42' = '5.73'
    'E42' = '  -2.63%  '
    'B43' = 'RenderToken'
    'C43' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D43' = '1.86'
    'E43' = '  +6.62%  '
    'B44' = 'PaxDollar'
    'C44' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    'D44' = '0.973'
    'E44' = '  -3.30%  '
    'D45' = '0.965'
    'E45' = '  -3.34%  '
    'D46' = '64.64'
    'E46' = '  +0.13%  '
    'D47' = '1.721.63'
    'E47' = '  +0.78%  '
    'D48' = '87.05'
    'E48' = '  +1.50%  '
    'D49' = '0.0525'
    'E49' = '  +1.38%  '
    'B50' = 'BitcoinSV'
    'C50' = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
    'D50' = '40.74'
    'E50' = '  +19.34%  '
    'B51' = 'BabyDogeCoin'
    'C51' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D51' = '0.0₇0980'
    'E51' = '  -4.55%  '
}

# Cells in columns D/E need to stay text; force a Text number format
# before writing so numeric-looking strings are not auto-converted.
$textRefs = @(
    'D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6',
    'D7', 'E7', 'D8', 'E8', 'E9', 'E10', 'D11', 'E11', 'D12', 'E12',
    'D13', 'E13', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18',
    'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23',
    'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28',
    'E28', 'D29', 'E29', 'E30', 'E31', 'E32', 'D33', 'E33', 'D34', 'E34',
    'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'D39', 'E39', 'D40', 'E40',
    'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45',
    'D46', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50',
    'D51', 'E51'
)

foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

foreach ($ref in $newValues.Keys) {
    $ws.Range($ref).Value = $newValues[$ref]
}

# Restore the original (style-less) formatting on the touched D/E cells
# by copying the style from an untouched plain-text cell (B2), so only
# the cell *values* differ from the starting workbook.
foreach ($ref in $textRefs) {
    $ws.Range($ref).Style = $ws.Range("B2").Style
}
